$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepend "Bireysel," to every existing "Rol" (role) value in column C.
$ws.Range("C2").Value = "Bireysel,Takım üyesi,Yazılımcı"
$ws.Range("C3").Value = "Bireysel,Yazılımcı"
$ws.Range("C4").Value = "Bireysel,Yardımcı"
$ws.Range("C5").Value = "Bireysel,Araştırmacı,Yardımcı"
$ws.Range("C6").Value = "Bireysel,Yazılımcı,Takım üyesi"
$ws.Range("C7").Value = "Bireysel,Yardımcı"
$ws.Range("C8").Value = "Bireysel,Yazılımcı,Yardımcı"
$ws.Range("C9").Value = "Bireysel,Takım üyesi,Araştırmacı"

# Re-normalize B9's cell formatting back to the column's plain/default style
# (it previously carried a stray distinct style index even though the
# rendered formatting is identical to the rest of the column).
$ws.Range("A9").Copy()
$ws.Range("B9").PasteSpecial(-4122)
